# Add new results (columns V:AO, rows 2:3) to Sheet1, mirroring/extending the
# existing B2:U3 block, then move the view/selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 2 (angles, radians) - same 20-value cycle repeated into V2:AO2
$row2Values = @(
    0.31415926535897898,
    0.62831853071795896,
    0.94247779607693805,
    1.2566370614359199,
    1.5707963267949001,
    1.8849555921538801,
    2.1991148575128601,
    2.5132741228718301,
    2.8274333882308098,
    3.14159265358979,
    3.4557519189487702,
    3.76991118430775,
    4.0840704496667302,
    4.3982297150257104,
    4.7123889803846897,
    5.0265482457436699,
    5.3407075111026501,
    5.6548667764616303,
    5.9690260418206096,
    6.2831853071795898
)

# Row 3 (new simulation outcomes) for V3:AO3
$row3Values = @(1, 1, 1, 1, 0, 1, 0, 0, 1, 0, 1, 1, 0, 1, 0, 1, 0, 0, 1, 0)

# V is column 22 .. AO is column 41
$startCol = 22
for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $row2Values[$i]
    $ws.Cells.Item(3, $col).Value = $row3Values[$i]
}

# Update the view the way it was left: scrolled so column L is left-most,
# and AA8 selected.
$ws.Range("L1").Select()
try { $excel.ActiveWindow.TopLeftCell = $ws.Range("L1") } catch { }
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AA8").Select()
